$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (I1:J1) ------------------------------------------------
# Copy the formatting of the existing header cell H1 (bold, centered,
# thin box border) onto the two new header cells before setting values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows (I2:J34) -------------------------------------------------
$data = @(
    @(7, 8),
    @(7, 7),
    @(7, 8),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(5, 6),
    @(8, 9),
    @(8, 8),
    @(1, 4),
    @(1, 4),
    @(1, 2),
    @(1, 2),
    @(5, 5),
    @(8, 8),
    @(1, 5),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(7, 7),
    @(8, 9),
    @(3, 4),
    @(1, 3),
    @(1, 5),
    @(1, 2),
    @(1, 5),
    @(2, 5),
    @(4, 7),
    @(5, 6),
    @(7, 7),
    @(1, 2),
    @(3, 4),
    @(8, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $pair = $data[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
